# Applies the "Cosmic Symphony / Astronomy" -> "Enchanting Realm / History"
# rewrite described by the supplied diff.

$d = $word.ActiveDocument

function Replace-InRange {
    param($range, [string]$oldText, [string]$newText)
    $r = $range.Duplicate
    $found = $r.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
    if (-not $found) {
        throw "Could not find text: $oldText"
    }
    $r.Text = $newText
    return $r
}

# ---------------------------------------------------------------------
# Paragraph 1: Title
# ---------------------------------------------------------------------
Replace-InRange $d.Paragraphs(1).Range `
    "A Cosmic Symphony: Unveiling the Harmony of the Universe" `
    "The Enchanting Realm of History: A Journey Through Time" | Out-Null

# ---------------------------------------------------------------------
# Paragraph 2: Author name
# ---------------------------------------------------------------------
Replace-InRange $d.Paragraphs(2).Range "Amelia Newman" "Emily Rose" | Out-Null

# ---------------------------------------------------------------------
# Paragraph 3: Email address
#   "amelianewman@gmail" + "." + "com"
#   -> "emily" + "." + "rose@schoolmail" + "." + "edu"
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs(3).Range
Replace-InRange $p3 "amelianewman@gmail" "emily" | Out-Null
$comRange = Replace-InRange $p3 "com" "rose@schoolmail"
$comRange.Collapse(0)
$comRange.InsertAfter(".")
$comRange.Collapse(0)
$comRange.InsertAfter("edu")

# ---------------------------------------------------------------------
# Paragraph 5: first body paragraph
# ---------------------------------------------------------------------
$p5 = $d.Paragraphs(5).Range

Replace-InRange $p5 `
    "Astronomy, an enthralling field that captivates the imagination, has unraveled the complexities of the cosmos" `
    "History, the captivating saga of humanity's journey through time, invites us to embark on an enthralling voyage of discovery" | Out-Null

Replace-InRange $p5 `
    " From the birth and evolution of stars to the expansion of the universe, the science of celestial bodies unveils a cosmic symphony brimming with intriguing melodies and astonishing rhythms" `
    " Each chapter of this grand narrative unfolds like a tapestry woven with tales of courage, resilience, and innovation" | Out-Null

$r = Replace-InRange $p5 `
    " Through telescopes and observatories, astronomers embark on a journey of discovery, unearthing celestial wonders that challenge our understanding of time, space, and existence" `
    " From the dawn of civilization to the modern era, history unveils the intricate interconnectedness of human experiences, shedding light on our origins, our struggles, and our triumphs"
$r.Collapse(0)
$r.InsertAfter(".")
$r.Collapse(0)
$r.InsertAfter(" It whispers of forgotten empires and lost civilizations, of epic battles and peaceful revolutions")
$r.Collapse(0)
$r.InsertAfter(".")
$r.Collapse(0)
$r.InsertAfter(" Within its vast embrace, history holds the keys to understanding our present and forging a better future")

Replace-InRange $p5 `
    "As we traverse the cosmic landscape, we encounter mesmerizing celestial bodies, each contributing a unique voice to the universal chorus" `
    "In the annals of history, we encounter iconic figures who shaped the course of events, leaving an indelible mark on the world" | Out-Null

Replace-InRange $p5 `
    " Planets, with their diverse characteristics, dance around their host stars, forming intricate choreographies" `
    " Their decisions, driven by ambition, compassion, or folly, reverberated across centuries, shaping societies and cultures" | Out-Null

Replace-InRange $p5 `
    " Stars, like fiery beacons, illuminate the darkness, emitting captivating radiations that paint ethereal tapestries across the night sky" `
    " We learn from their successes and failures, gaining insights into the complexities of human nature" | Out-Null

$r = Replace-InRange $p5 `
    " Galaxies, vast conglomerations of stars, gas, and dust, emerge as sprawling metropolises teeming with celestial activity" `
    " History also introduces us to ordinary individuals whose lives, though humble, collectively weave the rich fabric of the past"
$r.Collapse(0)
$r.InsertAfter(".")
$r.Collapse(0)
$r.InsertAfter(" Their stories, often overlooked, remind us that every person has a role to play in the grand scheme of things")

Replace-InRange $p5 `
    "The harmony of the universe extends beyond the visible realm" `
    "Furthermore, history teaches us about the intricate interplay between individuals and the forces that shape their lives" | Out-Null

Replace-InRange $p5 `
    " Electromagnetic waves, spanning a spectrum from radio waves to gamma rays, permeate the cosmos, carrying vital information about celestial objects" `
    " It reveals the impact of geography, climate, and technology on human societies" | Out-Null

Replace-InRange $p5 `
    " Gravity, an invisible force, orchestrates the celestial ballet, guiding the motion of planets, stars, and galaxies" `
    " By examining past events, we gain a deeper understanding of the challenges and opportunities that confront us today" | Out-Null

Replace-InRange $p5 `
    " Through these cosmic connections, the universe resonates with an underlying unity, revealing hidden patterns and profound interconnectedness" `
    " History equips us with the critical thinking skills necessary to navigate an ever-changing world, enabling us to make informed decisions and contribute meaningfully to society" | Out-Null

# ---------------------------------------------------------------------
# Paragraph 6: "Summary" heading -- unchanged
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# Paragraph 7: Summary body paragraph
# ---------------------------------------------------------------------
$p7 = $d.Paragraphs(7).Range

Replace-InRange $p7 `
    "Astronomy, with its enchanting allure, unveils the intricate symphony of the cosmos" `
    "History, an enthralling narrative of humanity's journey through time, unveils the intricate interconnectedness of human experiences" | Out-Null

Replace-InRange $p7 `
    " Gazing into the vast expanse of the universe, we discover an enchanting harmony woven from celestial bodies, cosmic phenomena, and physical forces" `
    " It invites us to learn from the past, gaining insights into our origins, our struggles, and our triumphs" | Out-Null

Replace-InRange $p7 `
    " From the celestial dance of planets to the mesmerizing radiance of stars, the universe unfolds a story of unity, complexity, and awe-inspiring beauty" `
    " Through the study of history, we encounter iconic figures and ordinary individuals whose actions shaped the course of events" | Out-Null

$r = Replace-InRange $p7 `
    " As we continue to explore and comprehend the cosmos, we unveil the interconnectedness of life and the profound interconnectedness of all things" `
    " We delve into "
$r.Collapse(0)
$r.InsertAfter("the interplay between individuals and the forces that shape their lives, gaining a deeper understanding of the challenges and opportunities that confront us today")
$r.Collapse(0)
$r.InsertAfter(".")
$r.Collapse(0)
$r.InsertAfter(" History equips us with critical thinking skills, enabling us to make informed decisions and contribute meaningfully to society")

# Trailing "." run at the end of paragraph 7 is unchanged.

# ---------------------------------------------------------------------
# New empty paragraph at the very end of the document
# ---------------------------------------------------------------------
$endRange = $d.Paragraphs($d.Paragraphs.Count).Range
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
